# Updated symbol list on Wed Dec 28 19:48:24 UTC 2022 with GitHub Actions
# Applies refreshed "Price" (column D) quotes, one "Volume(1h)" label fix
# (E18), and the BKEXToken / CEJI / KickToken row reshuffle (rows 41-43)
# to the crypto price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") cells are stored as TEXT (inline strings), not
# numbers, in the source workbook. Plain `.Value = "123.45"` assignment
# would let the engine auto-coerce a numeric-looking string into a real
# number, which would change the stored cell type. Force text by
# formatting the touched cells as Text ("@") first, then clear the
# formatting afterwards so no stray number-format is left behind on the
# cells (matches the original, unformatted look of the data rows).

$priceUpdates = [ordered]@{
    "D2"  = "243.66"
    "D3"  = "23.95"
    "D4"  = "5.243"
    "D5"  = "0.05817"
    "D7"  = "3.231"
    "D8"  = "0.8083"
    "D9"  = "0.8862"
    "D12" = "0.03055"
    "D13" = "0.03039"
    "D14" = "0.09332"
    "D15" = "3.830"
    "D16" = "0.001539"
    "D17" = "0.04718"
    "D19" = "0.006176"
    "D20" = "0.001258"
    "D21" = "0.004071"
    "D24" = "2.157"
    "D25" = "0.3183"
    "D26" = "0.1328"
    "D40" = "0.03861"
    "D41" = "0.006255"
    "D42" = "0.1051"
    "D43" = "0.002516"
    "D44" = "0.007836"
    "D45" = "0.00005332"
    "D47" = "0.5354"
    "D48" = "0.003086"
    "D50" = "0.0002001"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# --- Row 18: "Worstin24h" suffix added to the Volume(1h) label.
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Rows 41-43: the ranking reshuffled, rotating Coin/Link/Volume text
# among the three rows (BKEXToken -> row42, CEJI -> row43, KickToken ->
# row41).

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
